$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# IDEFICS-5: changed conversion number in tg, chol, hdl and ldl slightly to be consistent with EPICP
$ws.Range("H2").Value = "trg/88.57"
$ws.Range("H3").Value = "tc/38.67"
$ws.Range("H4").Value = "ldl/38.67"
$ws.Range("H5").Value = "hdl/38.67"

[void]$ws.Range("F31:F32").Select()
